$wb = $excel.ActiveWorkbook

# 1. Create the new "managers" sheet by copying "people" (places it right after "people")
$peopleWs = $wb.Worksheets.Item("people")
$peopleWs.Copy($null, $peopleWs) | Out-Null
$managersWs = $wb.Worksheets.Item("people (2)")
$managersWs.Name = "managers"

# Remove the extra rows copied from "people" (it only needs 4 data rows)
$managersWs.Range("A6:E9").Clear()

# Fill in the new manager records (written column-by-column to match authoring order)
$managersWs.Cells.Item(2,1).Value = 100
$managersWs.Cells.Item(3,1).Value = 101
$managersWs.Cells.Item(4,1).Value = 102
$managersWs.Cells.Item(5,1).Value = 103

$managersWs.Cells.Item(2,2).Value = "carol"
$managersWs.Cells.Item(3,2).Value = "brian"
$managersWs.Cells.Item(4,2).Value = "stephanie"
$managersWs.Cells.Item(5,2).Value = "dylan"

$managersWs.Cells.Item(2,3).Value = "alpha"
$managersWs.Cells.Item(3,3).Value = "bravo"
$managersWs.Cells.Item(4,3).Value = "charlie"
$managersWs.Cells.Item(5,3).Value = "delta"

$managersWs.Cells.Item(2,4).Value = 54
$managersWs.Cells.Item(3,4).Value = 64
$managersWs.Cells.Item(4,4).Value = 14
$managersWs.Cells.Item(5,4).Value = 21

$managersWs.Cells.Item(2,5).Value = "Copenhagen"
$managersWs.Cells.Item(3,5).Value = "Alger"
$managersWs.Cells.Item(4,5).Value = "Cairo"
$managersWs.Cells.Item(5,5).Value = "Lima"

$managersWs.Range("E6").Select() | Out-Null

# 2. Update the "expenses" sheet: tweak a value and move the selection there
$expensesWs = $wb.Worksheets.Item("expenses")
$expensesWs.Range("A11").Value = 102
$expensesWs.Activate() | Out-Null
$expensesWs.Range("A7").Select() | Out-Null
